$d = $word.ActiveDocument

$d.Content.Find.Execute("195×7=1365", $true, $false, $false, $false, $false, $true, 1, $false, "761×5=3805", 2) | Out-Null
$d.Content.Find.Execute("552×8=4416", $true, $false, $false, $false, $false, $true, 1, $false, "230×3=690", 2) | Out-Null
$d.Content.Find.Execute("709×9=6381", $true, $false, $false, $false, $false, $true, 1, $false, "367×5=1835", 2) | Out-Null
$d.Content.Find.Execute("248×6=1488", $true, $false, $false, $false, $false, $true, 1, $false, "844×9=7596", 2) | Out-Null
$d.Content.Find.Execute("234×7=1638", $true, $false, $false, $false, $false, $true, 1, $false, "418×2=836", 2) | Out-Null
$d.Content.Find.Execute("305×9=2745", $true, $false, $false, $false, $false, $true, 1, $false, "181×5=905", 2) | Out-Null
$d.Content.Find.Execute("260×4=1040", $true, $false, $false, $false, $false, $true, 1, $false, "261×3=783", 2) | Out-Null
$d.Content.Find.Execute("409×3=1227", $true, $false, $false, $false, $false, $true, 1, $false, "423×6=2538", 2) | Out-Null
$d.Content.Find.Execute("314×7=2198", $true, $false, $false, $false, $false, $true, 1, $false, "721×3=2163", 2) | Out-Null
$d.Content.Find.Execute("149×4=596", $true, $false, $false, $false, $false, $true, 1, $false, "935×2=1870", 2) | Out-Null
$d.Content.Find.Execute("243×5=1215", $true, $false, $false, $false, $false, $true, 1, $false, "698×2=1396", 2) | Out-Null
$d.Content.Find.Execute("386×2=772", $true, $false, $false, $false, $false, $true, 1, $false, "484×8=3872", 2) | Out-Null
$d.Content.Find.Execute("182×7=1274", $true, $false, $false, $false, $false, $true, 1, $false, "513×2=1026", 2) | Out-Null
$d.Content.Find.Execute("976×5=4880", $true, $false, $false, $false, $false, $true, 1, $false, "814×4=3256", 2) | Out-Null
$d.Content.Find.Execute("778×5=3890", $true, $false, $false, $false, $false, $true, 1, $false, "620×9=5580", 2) | Out-Null
$d.Content.Find.Execute("656×2=1312", $true, $false, $false, $false, $false, $true, 1, $false, "577×6=3462", 2) | Out-Null
$d.Content.Find.Execute("554×3=1662", $true, $false, $false, $false, $false, $true, 1, $false, "944×2=1888", 2) | Out-Null
$d.Content.Find.Execute("120×4=480", $true, $false, $false, $false, $false, $true, 1, $false, "609×5=3045", 2) | Out-Null
$d.Content.Find.Execute("589×7=4123", $true, $false, $false, $false, $false, $true, 1, $false, "144×3=432", 2) | Out-Null
$d.Content.Find.Execute("930×4=3720", $true, $false, $false, $false, $false, $true, 1, $false, "199×3=597", 2) | Out-Null
$d.Content.Find.Execute("856×2=1712", $true, $false, $false, $false, $false, $true, 1, $false, "612×5=3060", 2) | Out-Null
$d.Content.Find.Execute("457×4=1828", $true, $false, $false, $false, $false, $true, 1, $false, "138×6=828", 2) | Out-Null
$d.Content.Find.Execute("514×8=4112", $true, $false, $false, $false, $false, $true, 1, $false, "891×7=6237", 2) | Out-Null
$d.Content.Find.Execute("613×6=3678", $true, $false, $false, $false, $false, $true, 1, $false, "759×9=6831", 2) | Out-Null
$d.Content.Find.Execute("187×4=748", $true, $false, $false, $false, $false, $true, 1, $false, "239×6=1434", 2) | Out-Null
